$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append new row 27 with the new mail-log entry
$row = 27
$ws.Cells.Item($row, 1).Value = "Leg dit even neer bij Koen."
$ws.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value = "Testmail #15: Leg dit even neer bij Koen."
$ws.Cells.Item($row, 4).Value = "Planning / Afspraak"
$ws.Cells.Item($row, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$ws.Cells.Item($row, 6).Value = "2025-08-04 20:58:04"
$ws.Cells.Item($row, 7).Value = "Ja"
$ws.Cells.Item($row, 8).Value = "Ja"
$ws.Cells.Item($row, 9).Value = "Nee"
$ws.Cells.Item($row, 10).Value = "Nee"

# Extend the conditional-formatting ranges that covered rows 2:26 so they
# now cover the newly added row 27 as well.
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $ws.Range($col + "2:" + $col + "26")
    $newRange = $ws.Range($col + "2:" + $col + "27")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for "Planning / Afspraak" (6 -> 7)
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(2, 2).Value = 7
